$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp text in A1 ---
$ws.Range("A1").Value = "Datos actualizados a 5 de Abril de 2020 a las 23:22"

# --- Row 4: Estados Unidos - updated totals ---
$ws.Range("B4").Value = 333173
$ws.Range("C4").Value = 21816
$ws.Range("D4").Value = 17177
$ws.Range("E4").Value = 306460
$ws.Range("F4").Value = 8623
$ws.Range("G4").Value = 1084
$ws.Range("H4").Value = 9536

# --- Rows 102-104: Nigeria moved above Malta/Mauricio in country list,
#     with Nigeria receiving fresh data and Malta/Mauricio shifting down ---
$ws.Range("A102").Value = "Nigeria"
$ws.Range("B102").Value = 232
$ws.Range("C102").Value = 18
$ws.Range("D102").Value = 33
$ws.Range("E102").Value = 194
$ws.Range("F102").Value = 2
$ws.Range("G102").Value = 1
$ws.Range("H102").Value = 5

$ws.Range("A103").Value = "Malta"
$ws.Range("B103").Value = 227
$ws.Range("C103").Value = 14
$ws.Range("D103").Value = 5
$ws.Range("E103").Value = 222
$ws.Range("F103").Value = 3
$ws.Range("G103").Value = 0
$ws.Range("H103").Value = 0

$ws.Range("A104").Value = "Mauricio"
$ws.Range("B104").Value = 227
$ws.Range("C104").Value = 31
$ws.Range("D104").Value = 7
$ws.Range("E104").Value = 213
$ws.Range("F104").Value = 1
$ws.Range("G104").Value = 0
$ws.Range("H104").Value = 7

# --- Row 130: Monaco - F column updated ---
$ws.Range("F130").Value = 4

# --- Rows 151-152: San Martin (Parte Francesa) moved above Eritrea,
#     with San Martin receiving fresh data and Eritrea shifting down ---
$ws.Range("A151").Value = "San Martin (Parte Francesa)"
$ws.Range("B151").Value = 32
$ws.Range("C151").Value = 8
$ws.Range("D151").Value = 7
$ws.Range("E151").Value = 23
$ws.Range("F151").Value = 6
$ws.Range("G151").Value = 0
$ws.Range("H151").Value = 2

$ws.Range("A152").Value = "Eritrea"
$ws.Range("B152").Value = 29
$ws.Range("C152").Value = 0
$ws.Range("D152").Value = 0
$ws.Range("E152").Value = 29
$ws.Range("F152").Value = 0
$ws.Range("G152").Value = 0
$ws.Range("H152").Value = 0

# --- Rows 165-166: Nueva Caledonia moved above Islas Virgenes de los
#     Estados Unidos, with Nueva Caledonia receiving fresh data and
#     Islas Virgenes shifting down ---
$ws.Range("A165").Value = "Nueva Caledonia"
$ws.Range("B165").Value = 18
$ws.Range("C165").Value = 1
$ws.Range("D165").Value = 1
$ws.Range("E165").Value = 17
$ws.Range("F165").Value = 0
$ws.Range("G165").Value = 0
$ws.Range("H165").Value = 0

$ws.Range("A166").Value = "Islas Virgenes de los Estados Unidos"
$ws.Range("B166").Value = 17
$ws.Range("C166").Value = 0
$ws.Range("D166").Value = 0
$ws.Range("E166").Value = 17
$ws.Range("F166").Value = 0
$ws.Range("G166").Value = 0
$ws.Range("H166").Value = 0
